$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 95 / 96 got reordered upstream (match swapped on re-sort by kickoff
#    order) -- columns B..AC swap between the two rows, column A (row index)
#    stays put.
# ---------------------------------------------------------------------------
$ws.Cells.Item(95, 2).Value2  = 6782565
$ws.Range("F95").Value        = "Santos de Gupiles"
$ws.Range("G95").Value        = "Municipal Perez Zeledon"
$ws.Cells.Item(95, 8).Value2  = 2
$ws.Cells.Item(95, 9).Value2  = 0
$ws.Range("J95").Value        = "H"
$ws.Cells.Item(95, 11).Value2 = 1.833
$ws.Cells.Item(95, 12).Value2 = 3.4
$ws.Cells.Item(95, 13).Value2 = 3.75
$ws.Cells.Item(95, 14).Value2 = 1.833
$ws.Cells.Item(95, 15).Value2 = 3.5
$ws.Cells.Item(95, 16).Value2 = 3.75
$ws.Cells.Item(95, 17).Value2 = -0.5
$ws.Cells.Item(95, 18).Value2 = 1.875
$ws.Cells.Item(95, 19).Value2 = 1.925
$ws.Cells.Item(95, 20).Value2 = 2.75
$ws.Cells.Item(95, 21).Value2 = 2
$ws.Cells.Item(95, 22).Value2 = 1.8
$ws.Cells.Item(95, 23).Value2 = 0.833
$ws.Cells.Item(95, 24).Value2 = -1
$ws.Cells.Item(95, 25).Value2 = -1
$ws.Cells.Item(95, 26).Value2 = 0.875
$ws.Cells.Item(95, 27).Value2 = -1
$ws.Cells.Item(95, 28).Value2 = -1
$ws.Cells.Item(95, 29).Value2 = 0.8

$ws.Cells.Item(96, 2).Value2  = 6782567
$ws.Range("F96").Value        = "AD Grecia"
$ws.Range("G96").Value        = "Municipal Liberia"
$ws.Cells.Item(96, 8).Value2  = 2
$ws.Cells.Item(96, 9).Value2  = 3
$ws.Range("J96").Value        = "A"
$ws.Cells.Item(96, 11).Value2 = 2.875
$ws.Cells.Item(96, 12).Value2 = 3.5
$ws.Cells.Item(96, 13).Value2 = 2.15
$ws.Cells.Item(96, 14).Value2 = 2.3
$ws.Cells.Item(96, 15).Value2 = 3.5
$ws.Cells.Item(96, 16).Value2 = 2.6
$ws.Cells.Item(96, 17).Value2 = 0
$ws.Cells.Item(96, 18).Value2 = 1.8
$ws.Cells.Item(96, 19).Value2 = 2
$ws.Cells.Item(96, 20).Value2 = 2.75
$ws.Cells.Item(96, 21).Value2 = 1.8
$ws.Cells.Item(96, 22).Value2 = 2
$ws.Cells.Item(96, 23).Value2 = -1
$ws.Cells.Item(96, 24).Value2 = -1
$ws.Cells.Item(96, 25).Value2 = 1.6
$ws.Cells.Item(96, 26).Value2 = -1
$ws.Cells.Item(96, 27).Value2 = 1
$ws.Cells.Item(96, 28).Value2 = 0.8
$ws.Cells.Item(96, 29).Value2 = -1

# ---------------------------------------------------------------------------
# 2) Rows 110 / 111 -- same kind of swap.
# ---------------------------------------------------------------------------
$ws.Cells.Item(110, 2).Value2  = 6782581
$ws.Range("F110").Value        = "Alajuelense"
$ws.Range("G110").Value        = "AD Grecia"
$ws.Cells.Item(110, 8).Value2  = 2
$ws.Cells.Item(110, 9).Value2  = 0
$ws.Range("J110").Value        = "H"
$ws.Cells.Item(110, 11).Value2 = 1.181
$ws.Cells.Item(110, 12).Value2 = 6.5
$ws.Cells.Item(110, 13).Value2 = 11
$ws.Cells.Item(110, 14).Value2 = 1.25
$ws.Cells.Item(110, 15).Value2 = 5
$ws.Cells.Item(110, 16).Value2 = 9
$ws.Cells.Item(110, 17).Value2 = -1.75
$ws.Cells.Item(110, 18).Value2 = 1.975
$ws.Cells.Item(110, 19).Value2 = 1.825
$ws.Cells.Item(110, 20).Value2 = 3.25
$ws.Cells.Item(110, 21).Value2 = 2
$ws.Cells.Item(110, 22).Value2 = 1.8
$ws.Cells.Item(110, 23).Value2 = 0.25
$ws.Cells.Item(110, 24).Value2 = -1
$ws.Cells.Item(110, 25).Value2 = -1
$ws.Cells.Item(110, 26).Value2 = 0.4875
$ws.Cells.Item(110, 27).Value2 = -0.5
$ws.Cells.Item(110, 28).Value2 = -1
$ws.Cells.Item(110, 29).Value2 = 0.8

$ws.Cells.Item(111, 2).Value2  = 6782579
$ws.Range("F111").Value        = "Santos de Gupiles"
$ws.Range("G111").Value        = "AD San Carlos"
$ws.Cells.Item(111, 8).Value2  = 0
$ws.Cells.Item(111, 9).Value2  = 2
$ws.Range("J111").Value        = "A"
$ws.Cells.Item(111, 11).Value2 = 2.4
$ws.Cells.Item(111, 12).Value2 = 3.3
$ws.Cells.Item(111, 13).Value2 = 2.7
$ws.Cells.Item(111, 14).Value2 = 2.375
$ws.Cells.Item(111, 15).Value2 = 3.4
$ws.Cells.Item(111, 16).Value2 = 2.8
$ws.Cells.Item(111, 17).Value2 = -0.25
$ws.Cells.Item(111, 18).Value2 = 2
$ws.Cells.Item(111, 19).Value2 = 1.8
$ws.Cells.Item(111, 20).Value2 = 2.5
$ws.Cells.Item(111, 21).Value2 = 1.875
$ws.Cells.Item(111, 22).Value2 = 1.925
$ws.Cells.Item(111, 23).Value2 = -1
$ws.Cells.Item(111, 24).Value2 = -1
$ws.Cells.Item(111, 25).Value2 = 1.8
$ws.Cells.Item(111, 26).Value2 = -1
$ws.Cells.Item(111, 27).Value2 = 0.8
$ws.Cells.Item(111, 28).Value2 = -1
$ws.Cells.Item(111, 29).Value2 = 0.925

# ---------------------------------------------------------------------------
# 3) Rows 192 / 193 -- same kind of swap.
# ---------------------------------------------------------------------------
$ws.Cells.Item(192, 2).Value2  = 7623919
$ws.Range("F192").Value        = "Municipal Liberia"
$ws.Range("G192").Value        = "Sporting San Jose"
$ws.Cells.Item(192, 8).Value2  = 2
$ws.Cells.Item(192, 9).Value2  = 0
$ws.Range("J192").Value        = "H"
$ws.Cells.Item(192, 11).Value2 = 1.75
$ws.Cells.Item(192, 12).Value2 = 3.6
$ws.Cells.Item(192, 13).Value2 = 3.8
$ws.Cells.Item(192, 14).Value2 = 1.8
$ws.Cells.Item(192, 15).Value2 = 3.6
$ws.Cells.Item(192, 16).Value2 = 3.6
$ws.Cells.Item(192, 17).Value2 = -0.5
$ws.Cells.Item(192, 18).Value2 = 1.9
$ws.Cells.Item(192, 19).Value2 = 1.9
$ws.Cells.Item(192, 20).Value2 = 2.75
$ws.Cells.Item(192, 21).Value2 = 2
$ws.Cells.Item(192, 22).Value2 = 1.8
$ws.Cells.Item(192, 23).Value2 = 0.8
$ws.Cells.Item(192, 24).Value2 = -1
$ws.Cells.Item(192, 25).Value2 = -1
$ws.Cells.Item(192, 26).Value2 = 0.8999999999999999
$ws.Cells.Item(192, 27).Value2 = -1
$ws.Cells.Item(192, 28).Value2 = -1
$ws.Cells.Item(192, 29).Value2 = 0.8

$ws.Cells.Item(193, 2).Value2  = 7623916
$ws.Range("F193").Value        = "Santos de Gupiles"
$ws.Range("G193").Value        = "AD Grecia"
$ws.Cells.Item(193, 8).Value2  = 0
$ws.Cells.Item(193, 9).Value2  = 2
$ws.Range("J193").Value        = "A"
$ws.Cells.Item(193, 11).Value2 = 2.05
$ws.Cells.Item(193, 12).Value2 = 3.3
$ws.Cells.Item(193, 13).Value2 = 3.2
$ws.Cells.Item(193, 14).Value2 = 1.909
$ws.Cells.Item(193, 15).Value2 = 3.4
$ws.Cells.Item(193, 16).Value2 = 3.6
$ws.Cells.Item(193, 17).Value2 = -0.5
$ws.Cells.Item(193, 18).Value2 = 1.95
$ws.Cells.Item(193, 19).Value2 = 1.85
$ws.Cells.Item(193, 20).Value2 = 2.5
$ws.Cells.Item(193, 21).Value2 = 1.85
$ws.Cells.Item(193, 22).Value2 = 1.95
$ws.Cells.Item(193, 23).Value2 = -1
$ws.Cells.Item(193, 24).Value2 = -1
$ws.Cells.Item(193, 25).Value2 = 2.6
$ws.Cells.Item(193, 26).Value2 = -1
$ws.Cells.Item(193, 27).Value2 = 0.8500000000000001
$ws.Cells.Item(193, 28).Value2 = -1
$ws.Cells.Item(193, 29).Value2 = 0.95

# ---------------------------------------------------------------------------
# 4) Four newly-recorded matches appended as rows 195..198. Clone the
#    formatting (bold/bordered id column, date-formatted date column) from
#    the last existing row before filling in the values.
# ---------------------------------------------------------------------------
$ws.Range("A194:AC194").Copy()
$ws.Range("A195:AC198").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 195
$ws.Cells.Item(195, 1).Value2  = 193
$ws.Cells.Item(195, 2).Value2  = 7623931
$ws.Range("C195").Value        = "Costa Rica Primera Division"
$ws.Range("D195").Value        = "Costa Rica Primera Division"
$ws.Cells.Item(195, 5).Value2  = 45349.91666666666
$ws.Range("F195").Value        = "AD San Carlos"
$ws.Range("G195").Value        = "Santos de Gupiles"
$ws.Cells.Item(195, 8).Value2  = 4
$ws.Cells.Item(195, 9).Value2  = 0
$ws.Range("J195").Value        = "H"
$ws.Cells.Item(195, 11).Value2 = 1.333
$ws.Cells.Item(195, 12).Value2 = 4.5
$ws.Cells.Item(195, 13).Value2 = 7
$ws.Cells.Item(195, 14).Value2 = 1.363
$ws.Cells.Item(195, 15).Value2 = 4.5
$ws.Cells.Item(195, 16).Value2 = 6.5
$ws.Cells.Item(195, 17).Value2 = -1.25
$ws.Cells.Item(195, 18).Value2 = 1.875
$ws.Cells.Item(195, 19).Value2 = 1.925
$ws.Cells.Item(195, 20).Value2 = 2.75
$ws.Cells.Item(195, 21).Value2 = 1.9
$ws.Cells.Item(195, 22).Value2 = 1.9
$ws.Cells.Item(195, 23).Value2 = 0.363
$ws.Cells.Item(195, 24).Value2 = -1
$ws.Cells.Item(195, 25).Value2 = -1
$ws.Cells.Item(195, 26).Value2 = 0.875
$ws.Cells.Item(195, 27).Value2 = -1
$ws.Cells.Item(195, 28).Value2 = 0.8999999999999999
$ws.Cells.Item(195, 29).Value2 = -1

# Row 196
$ws.Cells.Item(196, 1).Value2  = 194
$ws.Cells.Item(196, 2).Value2  = 7623932
$ws.Range("C196").Value        = "Costa Rica Primera Division"
$ws.Range("D196").Value        = "Costa Rica Primera Division"
$ws.Cells.Item(196, 5).Value2  = 45350.75
$ws.Range("F196").Value        = "AD Guanacasteca"
$ws.Range("G196").Value        = "Cartagines"
$ws.Cells.Item(196, 8).Value2  = 1
$ws.Cells.Item(196, 9).Value2  = 0
$ws.Range("J196").Value        = "H"
$ws.Cells.Item(196, 11).Value2 = 2.5
$ws.Cells.Item(196, 12).Value2 = 3.2
$ws.Cells.Item(196, 13).Value2 = 2.5
$ws.Cells.Item(196, 14).Value2 = 2.7
$ws.Cells.Item(196, 15).Value2 = 3.1
$ws.Cells.Item(196, 16).Value2 = 2.375
$ws.Cells.Item(196, 17).Value2 = 0
$ws.Cells.Item(196, 18).Value2 = 2.025
$ws.Cells.Item(196, 19).Value2 = 1.775
$ws.Cells.Item(196, 20).Value2 = 2.25
$ws.Cells.Item(196, 21).Value2 = 1.8
$ws.Cells.Item(196, 22).Value2 = 2
$ws.Cells.Item(196, 23).Value2 = 1.7
$ws.Cells.Item(196, 24).Value2 = -1
$ws.Cells.Item(196, 25).Value2 = -1
$ws.Cells.Item(196, 26).Value2 = 1.025
$ws.Cells.Item(196, 27).Value2 = -1
$ws.Cells.Item(196, 28).Value2 = -1
$ws.Cells.Item(196, 29).Value2 = 1

# Row 197
$ws.Cells.Item(197, 1).Value2  = 195
$ws.Cells.Item(197, 2).Value2  = 7623935
$ws.Range("C197").Value        = "Costa Rica Primera Division"
$ws.Range("D197").Value        = "Costa Rica Primera Division"
$ws.Cells.Item(197, 5).Value2  = 45350.89583333334
$ws.Range("F197").Value        = "Sporting San Jose"
$ws.Range("G197").Value        = "AD Grecia"
$ws.Cells.Item(197, 8).Value2  = 4
$ws.Cells.Item(197, 9).Value2  = 1
$ws.Range("J197").Value        = "H"
$ws.Cells.Item(197, 11).Value2 = 2.2
$ws.Cells.Item(197, 12).Value2 = 3
$ws.Cells.Item(197, 13).Value2 = 3.2
$ws.Cells.Item(197, 14).Value2 = 1.85
$ws.Cells.Item(197, 15).Value2 = 3.3
$ws.Cells.Item(197, 16).Value2 = 3.75
$ws.Cells.Item(197, 17).Value2 = -0.5
$ws.Cells.Item(197, 18).Value2 = 1.925
$ws.Cells.Item(197, 19).Value2 = 1.875
$ws.Cells.Item(197, 20).Value2 = 2.25
$ws.Cells.Item(197, 21).Value2 = 2
$ws.Cells.Item(197, 22).Value2 = 1.8
$ws.Cells.Item(197, 23).Value2 = 0.8500000000000001
$ws.Cells.Item(197, 24).Value2 = -1
$ws.Cells.Item(197, 25).Value2 = -1
$ws.Cells.Item(197, 26).Value2 = 0.925
$ws.Cells.Item(197, 27).Value2 = -1
$ws.Cells.Item(197, 28).Value2 = 1
$ws.Cells.Item(197, 29).Value2 = -1

# Row 198
$ws.Cells.Item(198, 1).Value2  = 196
$ws.Cells.Item(198, 2).Value2  = 7623933
$ws.Range("C198").Value        = "Costa Rica Primera Division"
$ws.Range("D198").Value        = "Costa Rica Primera Division"
$ws.Cells.Item(198, 5).Value2  = 45350.95833333334
$ws.Range("F198").Value        = "Herediano"
$ws.Range("G198").Value        = "Municipal Liberia"
$ws.Cells.Item(198, 8).Value2  = 4
$ws.Cells.Item(198, 9).Value2  = 0
$ws.Range("J198").Value        = "H"
$ws.Cells.Item(198, 11).Value2 = 1.533
$ws.Cells.Item(198, 12).Value2 = 4
$ws.Cells.Item(198, 13).Value2 = 5
$ws.Cells.Item(198, 14).Value2 = 1.5
$ws.Cells.Item(198, 15).Value2 = 4.2
$ws.Cells.Item(198, 16).Value2 = 5
$ws.Cells.Item(198, 17).Value2 = -1
$ws.Cells.Item(198, 18).Value2 = 1.95
$ws.Cells.Item(198, 19).Value2 = 1.85
$ws.Cells.Item(198, 20).Value2 = 2.5
$ws.Cells.Item(198, 21).Value2 = 1.875
$ws.Cells.Item(198, 22).Value2 = 1.925
$ws.Cells.Item(198, 23).Value2 = 0.5
$ws.Cells.Item(198, 24).Value2 = -1
$ws.Cells.Item(198, 25).Value2 = -1
$ws.Cells.Item(198, 26).Value2 = 0.95
$ws.Cells.Item(198, 27).Value2 = -1
$ws.Cells.Item(198, 28).Value2 = 0.875
$ws.Cells.Item(198, 29).Value2 = -1
